$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D11 becomes a formula referencing C11 instead of a hard-coded value
$ws.Range("D11").Formula = "=C11+2558"

# New row 18: another "Кормушка" / "для зяблика" entry, continuing the IP range
# right after row 15 (256..510), linked via formulas to the previous rows.
$ws.Range("A18").Value = "Кормушка"
$ws.Range("B18").Value = "для зяблика"
$ws.Range("C18").Formula = "=D15+1"
$ws.Range("D18").Formula = "=C18+254"
$ws.Range("E18").Value = "10.5.11.1"
$ws.Range("F18").Value = "10.5.11.255"
$ws.Range("G18").Value = "open"

# Match the saved selection shown in the diff
$ws.Range("D19").Select()
